$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.230.37"
$ws.Range("E2").Value = "  -1.11%  "
$ws.Range("D3").Value = "3.535.93"
$ws.Range("E3").Value = "  +0.53%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "607.83"
$ws.Range("E5").Value = "  -0.13%  "
$ws.Range("D6").Value = "143.74"
$ws.Range("E6").Value = "  -2.86%  "
$ws.Range("D7").Value = "3.535.96"
$ws.Range("E7").Value = "  +0.54%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +0.48%  "
$ws.Range("D10").Value = "0.137"
$ws.Range("E10").Value = "  -4.24%  "
$ws.Range("D11").Value = "8.06"
$ws.Range("E11").Value = "  +1.42%  "
$ws.Range("E12").Value = "  -2.69%  "
$ws.Range("D13").Value = "4.134.60"
$ws.Range("E13").Value = "  +0.57%  "
$ws.Range("E14").Value = "  -4.94%  "
$ws.Range("D15").Value = "30.31"
$ws.Range("D16").Value = "3.538.13"
$ws.Range("E16").Value = "  +0.76%  "
$ws.Range("D17").Value = "66.318.90"
$ws.Range("E17").Value = "  -1.03%  "
$ws.Range("E18").Value = "  -0.51%  "
$ws.Range("D19").Value = "10.90"
$ws.Range("E19").Value = "  +1.96%  "
$ws.Range("E20").Value = "  -3.60%  "
$ws.Range("D21").Value = "14.97"
$ws.Range("E21").Value = "  -2.34%  "
$ws.Range("D22").Value = "426.12"
$ws.Range("E22").Value = "  -2.71%  "
$ws.Range("D23").Value = "0.602"
$ws.Range("E23").Value = "  -1.20%  "
$ws.Range("D24").Value = "78.74"
$ws.Range("E24").Value = "  -0.72%  "
$ws.Range("D25").Value = "3.673.77"
$ws.Range("E25").Value = "  +0.48%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").Value = "0.0000120"
$ws.Range("E27").Value = "  -1.31%  "
$ws.Range("D28").Value = "8.07"
$ws.Range("E28").Value = "  -2.80%  "
$ws.Range("D29").Value = "9.20"
$ws.Range("E29").Value = "  -5.81%  "
$ws.Range("D30").Value = "2.47"
$ws.Range("E30").Value = "  -1.65%  "
$ws.Range("D31").Value = "0.998"
$ws.Range("E31").Value = "  -0.27%  "
$ws.Range("E32").Value = "  -3.67%  "
$ws.Range("E33").Value = "  -6.25%  "
$ws.Range("D34").Value = "25.31"
$ws.Range("E34").Value = "  -0.59%  "
$ws.Range("D35").Value = "3.524.80"
$ws.Range("E35").Value = "  +0.50%  "
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("D37").Value = "1.76"
$ws.Range("E37").Value = "  -2.79%  "
$ws.Range("E38").Value = "  -2.26%  "
$ws.Range("D39").Value = "5.62"
$ws.Range("E39").Value = "  -5.99%  "
$ws.Range("D41").Value = "172.33"
$ws.Range("E41").Value = "  -0.55%  "
$ws.Range("D42").Value = "0.0857"
$ws.Range("E42").Value = "  -4.12%  "
$ws.Range("D43").Value = "5.19"
$ws.Range("E43").Value = "  -4.38%  "
$ws.Range("E44").Value = "  -0.49%  "
$ws.Range("E45").Value = "  -8.81%  "
$ws.Range("D46").Value = "45.37"
$ws.Range("E46").Value = "  -1.51%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Value = "26.03"
$ws.Range("E47").Value = "  -7.19%  "
$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D48").Value = "1.21"
$ws.Range("E48").Value = "  -6.57%  "
$ws.Range("E49").Value = "  -2.37%  "
$ws.Range("E50").Value = "  -4.11%  "
$ws.Range("E51").Value = "  -4.64%  "
